$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels (shared strings)
$ws.Range("A1").Value = "legenda"
$ws.Range("B1").Value = "area"
$ws.Range("D1").Value = "area_km2"

# Update the B2 data value
$ws.Range("B2").Value = 174402.673922
